# IRPStudents.xlsx - UI clean-up pass
#
# The web app now only renders four columns (Status, Field, Expected Value,
# Actual Value) plus an extra "error" row under any field that failed
# validation. The STUDENT reference sheet's Excel view state had drifted
# (zoomed in at 150%, a stray selection left on AD1, and the Excel app
# window position) while someone was eyeballing column C
# (ResponsibleInstitutionIdentifier) - widen that column so its header is
# fully visible, and put the view back to a clean default state before
# saving.

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("STUDENT")
$win = $excel.ActiveWindow

[void]$ws.Activate()

# --- Application / workbook window position -----------------------------
# Nudge the app window back toward the left edge of the screen.
$win.Left = 20
$win.Top  = 0

# --- Sheet view: drop the 150% zoom back to the normal 100% -------------
$win.Zoom = 100

# --- Sheet view: clear the lingering selection, back to the home cell ---
[void]$ws.Range("A1").Select()

# --- Column widths --------------------------------------------------------
# Column C (ResponsibleInstitutionIdentifier) used to share the generic
# 10.83-wide default with columns D:I. Give it its own best-fit-sized
# column (~27 chars) so the full header text is visible, same treatment
# as the other bestFit columns (B, W, X, Z, AA, AD).
$ws.Columns.Item(3).ColumnWidth = 26.16666666666667
